$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...stránka, která se vrátí uživatelovi. " -> "...vrátí uživatelovi."
#    (drop the trailing space) and relocate the (singleton, hidden) _GoBack
#    bookmark to sit right after the new final period instead of in the
#    empty paragraph at the very end of the document.
# ---------------------------------------------------------------------------

# Locate the exact text "uživatelovi." with a search-only Find (no replace,
# so it can't disturb any later bookmark placement).
$loc = $d.Content
$loc.Find.Execute("uživatelovi.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$loc.Collapse(0)

# Re-adding a bookmark named "_GoBack" moves Word's singleton hidden
# "last edit" bookmark here (and removes it from wherever it used to be).
$d.Bookmarks.Add("_GoBack", $loc)

# Now remove the single trailing space character that used to follow the
# period (it now immediately follows the just-inserted bookmark).
$spacePos = $loc.Start
$spaceRange = $d.Range($spacePos, $spacePos + 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Text = ""
}

# ---------------------------------------------------------------------------
# 2) "...pro získání dat." -> "...pro získání dat (např. samostatný projekt
#    s databází)."
# ---------------------------------------------------------------------------

$loc2 = $d.Content
$loc2.Find.Execute("pro získání dat.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $loc2.Start
$end2 = $loc2.End
$r2 = $d.Range($start2, $end2)
$r2.Text = "pro získání dat (např. samostatný projekt s databází)."
